$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio6")

# Copy column K (G5:K35 pattern columns) formatting+values into L, M, N, O
# to replicate the same cell styles used for J/K columns.
$ws.Range("K5:K35").Copy($ws.Range("L5:L35"))
$ws.Range("K5:K35").Copy($ws.Range("M5:M35"))
$ws.Range("K5:K35").Copy($ws.Range("N5:N35"))
$ws.Range("K5:K35").Copy($ws.Range("O5:O35"))

# Match column width (21) used by columns H:K
$ws.Columns("L").ColumnWidth = 20 + 1/6
$ws.Columns("M").ColumnWidth = 20 + 1/6
$ws.Columns("N").ColumnWidth = 20 + 1/6
$ws.Columns("O").ColumnWidth = 20 + 1/6

# Fill in the generated values for the new road-delay configuration columns
$ws.Range("L6").Value = 90
$ws.Range("M6").Value = 93
$ws.Range("N6").Value = 88
$ws.Range("O6").Value = 91
$ws.Range("L7").Value = 78
$ws.Range("M7").Value = 85
$ws.Range("N7").Value = 81
$ws.Range("O7").Value = 80
$ws.Range("L8").Value = 75
$ws.Range("M8").Value = 78
$ws.Range("N8").Value = 84
$ws.Range("O8").Value = 78
$ws.Range("L9").Value = 72
$ws.Range("M9").Value = 81
$ws.Range("N9").Value = 77
$ws.Range("O9").Value = 77
$ws.Range("L10").Value = 72
$ws.Range("M10").Value = 84
$ws.Range("N10").Value = 74
$ws.Range("O10").Value = 76
$ws.Range("L11").Value = 75
$ws.Range("M11").Value = 78
$ws.Range("N11").Value = 76
$ws.Range("O11").Value = 73
$ws.Range("L12").Value = 76
$ws.Range("M12").Value = 76
$ws.Range("N12").Value = 74
$ws.Range("O12").Value = 71
$ws.Range("L13").Value = 79
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 84
$ws.Range("O13").Value = 78
$ws.Range("L14").Value = 76
$ws.Range("M14").Value = 81
$ws.Range("N14").Value = 79
$ws.Range("O14").Value = 77
$ws.Range("L15").Value = 81
$ws.Range("M15").Value = 78
$ws.Range("N15").Value = 74
$ws.Range("O15").Value = 73
$ws.Range("L16").Value = 77
$ws.Range("M16").Value = 76
$ws.Range("N16").Value = 76
$ws.Range("O16").Value = 79
$ws.Range("L17").Value = 72
$ws.Range("M17").Value = 83
$ws.Range("N17").Value = 81
$ws.Range("O17").Value = 71
$ws.Range("L18").Value = 80
$ws.Range("M18").Value = 71
$ws.Range("N18").Value = 72
$ws.Range("O18").Value = 78
$ws.Range("L19").Value = 74
$ws.Range("M19").Value = 83
$ws.Range("N19").Value = 80
$ws.Range("O19").Value = 76
$ws.Range("L20").Value = 75
$ws.Range("M20").Value = 72
$ws.Range("N20").Value = 70
$ws.Range("O20").Value = 75
$ws.Range("L21").Value = 73
$ws.Range("M21").Value = 74
$ws.Range("N21").Value = 81
$ws.Range("O21").Value = 72
$ws.Range("L22").Value = 79
$ws.Range("M22").Value = 80
$ws.Range("N22").Value = 75
$ws.Range("O22").Value = 83
$ws.Range("L23").Value = 74
$ws.Range("M23").Value = 86
$ws.Range("N23").Value = 78
$ws.Range("O23").Value = 78
$ws.Range("L24").Value = 77
$ws.Range("M24").Value = 78
$ws.Range("N24").Value = 77
$ws.Range("O24").Value = 77
$ws.Range("L25").Value = 74
$ws.Range("M25").Value = 76
$ws.Range("N25").Value = 86
$ws.Range("O25").Value = 76
$ws.Range("L26").Value = 69
$ws.Range("M26").Value = 74
$ws.Range("N26").Value = 72
$ws.Range("O26").Value = 76
$ws.Range("L27").Value = 74
$ws.Range("M27").Value = 77
$ws.Range("N27").Value = 77
$ws.Range("O27").Value = 78
$ws.Range("L28").Value = 77
$ws.Range("M28").Value = 83
$ws.Range("N28").Value = 80
$ws.Range("O28").Value = 80
$ws.Range("L29").Value = 78
$ws.Range("M29").Value = 78
$ws.Range("N29").Value = 73
$ws.Range("O29").Value = 73
$ws.Range("L30").Value = 82
$ws.Range("M30").Value = 81
$ws.Range("N30").Value = 73
$ws.Range("O30").Value = 73
$ws.Range("L31").Value = 90
$ws.Range("M31").Value = 93
$ws.Range("N31").Value = 88
$ws.Range("O31").Value = 91
$ws.Range("L32").Value = 23
$ws.Range("M32").Value = 24
$ws.Range("N32").Value = 24
$ws.Range("O32").Value = 16
$ws.Range("L33").Value = 10
$ws.Range("M33").Value = 12
$ws.Range("N33").Value = 11
$ws.Range("O33").Value = 10
$ws.Range("L34").Value = 57
$ws.Range("M34").Value = 57
$ws.Range("N34").Value = 53
$ws.Range("O34").Value = 65

# Row 35 holds the execution-time strings for each run
$ws.Range("L35").Value = "00:00:00:53272"
$ws.Range("M35").Value = "00:00:00:48693"
$ws.Range("N35").Value = "00:00:00:51418"
$ws.Range("O35").Value = "00:00:00:54981"

Write-Host "Road-delay generation columns (L:O) populated"
